$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.777.23'
$ws.Range('E2').Value = '  -0.74%  '

$ws.Range('D3').Value = '3.390.66'
$ws.Range('E3').Value = '  -3.00%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '577.71'

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.07'
$ws.Range('E6').Value = '  -5.61%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').Value = '3.390.53'
$ws.Range('E8').Value = '  -2.98%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -2.57%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.119'
$ws.Range('E10').Value = '  -9.81%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.00'
$ws.Range('E11').Value = '  -9.85%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.369'
$ws.Range('E12').Value = '  -7.69%  '

$ws.Range('D13').Value = '3.973.23'
$ws.Range('E13').Value = '  -2.88%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000174'
$ws.Range('E14').Value = '  -10.92%  '

$ws.Range('E15').Value = '  -1.71%  '

$ws.Range('D16').Value = '3.398.25'
$ws.Range('E16').Value = '  -3.35%  '

$ws.Range('D17').Value = '64.801.06'
$ws.Range('E17').Value = '  -0.76%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '25.80'
$ws.Range('E18').Value = '  -8.45%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.45'
$ws.Range('E19').Value = '  -13.94%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.73'
$ws.Range('E20').Value = '  -6.69%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.29'
$ws.Range('E21').Value = '  -6.17%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '376.08'
$ws.Range('E22').Value = '  -8.76%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.543'
$ws.Range('E23').Value = '  -8.23%  '

$ws.Range('E24').Value = '  -0.02%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '71.49'
$ws.Range('E25').Value = '  -7.28%  '

$ws.Range('D26').Value = '3.527.83'
$ws.Range('E26').Value = '  -3.05%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000102'
$ws.Range('E27').Value = '  -9.35%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.51%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.17'
$ws.Range('E29').Value = '  -10.40%  '

$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.90'
$ws.Range('E30').Value = '  -9.24%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.88'
$ws.Range('E31').Value = '  -9.60%  '

$ws.Range('D32').Value = '3.403.75'
$ws.Range('E32').Value = '  -2.79%  '

$ws.Range('E33').Value = '  +0.02%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '22.79'
$ws.Range('E34').Value = '  -5.37%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.140'
$ws.Range('E35').Value = '  -7.59%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '169.72'
$ws.Range('E36').Value = '  -1.75%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.57'
$ws.Range('E37').Value = '  -11.54%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.12'
$ws.Range('E38').Value = '  -11.50%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.42'
$ws.Range('E39').Value = '  -7.90%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.56'
$ws.Range('E40').Value = '  -12.46%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0741'
$ws.Range('E41').Value = '  -8.34%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.803'
$ws.Range('E42').Value = '  -5.38%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.09%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '42.34'
$ws.Range('E44').Value = '  -6.05%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.27'
$ws.Range('E45').Value = '  -14.89%  '

$ws.Range('E46').Value = '  -10.68%  '

$ws.Range('E47').Value = '  +1.06%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '21.63'
$ws.Range('E48').Value = '  -4.89%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.39'
$ws.Range('E49').Value = '  -8.23%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.141.57'
$ws.Range('E50').Value = '  -8.39%  '

$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.98'
$ws.Range('E51').Value = '  -15.28%  '
